# Insert a new data row at row 241 (pushes the existing rows 241-307 down
# to 242-308, growing the used range from A1:R307 to A1:R308).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(241).Insert()

# Populate the newly inserted row 241. It keeps the same market / product
# metadata as the row that used to occupy 241 (now shifted to 242), but it
# is a new weekly record with its own date and volume.
$ws.Cells.Item(241, 1).Value = 4
$ws.Cells.Item(241, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(241, 3).Value = "Los Lagos"
$ws.Cells.Item(241, 4).Value = 44736
$ws.Cells.Item(241, 5).Value = 10
$ws.Cells.Item(241, 6).Value = 100112040
$ws.Cells.Item(241, 7).Value = "Cilantro"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 200
$ws.Cells.Item(241, 11).Value = 12000
$ws.Cells.Item(241, 12).Value = 12000
$ws.Cells.Item(241, 13).Value = 12000
$ws.Cells.Item(241, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(241, 15).Value = "Región Metropolitana"
$ws.Cells.Item(241, 16).Value = 333
$ws.Cells.Item(241, 17).Value = 36
$ws.Cells.Item(241, 18).Value = "Hortaliza"
